$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the trailing "EDG 4930" run from the title paragraph.
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("EDG 4930", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = ""
}

# ------------------------------------------------------------------
# 2. Remove the "Scott Cohen" and "Charlotte Bolch" paragraphs
#    (the latter carried the _GoBack bookmark, which moves to the
#    due-dates table below).
# ------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Scott Cohen", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $scottPara = $r2.Paragraphs(1)
    $charlottePara = $scottPara.Next()
    $delRange = $d.Range($scottPara.Range.Start, $charlottePara.Range.End)
    $delRange.Delete()
}

# ------------------------------------------------------------------
# 3. Clear the due-date cells in the "Course Schedule & Assignment
#    Due Dates" table (the second table). The first cleared cell
#    keeps the _GoBack bookmark (now empty) in place of its text.
# ------------------------------------------------------------------
$t2 = $d.Tables(2)

# Row 2, column 3: add bookmark around the date text, then clear it.
$cell = $t2.Cell(2, 3)
$cellRange = $cell.Range
$textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
$d.Bookmarks.Add("_GoBack", $textRange) | Out-Null

$cell = $t2.Cell(2, 3)
$cellRange = $cell.Range
$textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
$textRange.Text = ""

# Rows 3-8, column 3: simply clear the date text.
for ($i = 3; $i -le 8; $i++) {
    $cell = $t2.Cell($i, 3)
    $cellRange = $cell.Range
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $textRange.Text = ""
}
